$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = -12.998
$ws.Range("B7").Value = 5.811000000000001
$ws.Range("D7").Value = -7.698
$ws.Range("A9").Value = -21.916
$ws.Range("D10").Value = -8.415000000000001
$ws.Range("B12").Value = 5.662
$ws.Range("A13").Value = -22.117
$ws.Range("D13").Value = -8.014999999999999
$ws.Range("B14").Value = 6.044
$ws.Range("C15").Value = -12.909
$ws.Range("A16").Value = -21.567
$ws.Range("D16").Value = -8.535
$ws.Range("A18").Value = -21.962
$ws.Range("B19").Value = 7.859999999999999
$ws.Range("A20").Value = -21.296
$ws.Range("D20").Value = -8.059999999999999
$ws.Range("D24").Value = -7.351999999999999
$ws.Range("A26").Value = -21.393
$ws.Range("B26").Value = 6.089
$ws.Range("A27").Value = -21.162
$ws.Range("B27").Value = 6.193000000000001
$ws.Range("C28").Value = -12.965
$ws.Range("A29").Value = -21.84
$ws.Range("B29").Value = 6.164000000000001
$ws.Range("D32").Value = -8.408000000000001
$ws.Range("C33").Value = -11.314
$ws.Range("A35").Value = -20.336
$ws.Range("C35").Value = -12.904
$ws.Range("A36").Value = -20.814
$ws.Range("B37").Value = 7.515000000000001
$ws.Range("B38").Value = 5.914
$ws.Range("C38").Value = -12.56
$ws.Range("D39").Value = -7.555
$ws.Range("C43").Value = -12.494
$ws.Range("C44").Value = -12.741
$ws.Range("A45").Value = -21.575
$ws.Range("C45").Value = -12.53
$ws.Range("B47").Value = 6.433
$ws.Range("C47").Value = -12.199
$ws.Range("D47").Value = -7.312
$ws.Range("D48").Value = -7.229000000000001
$ws.Range("B51").Value = 5.906
$ws.Range("C51").Value = -11.915
$ws.Range("B52").Value = 5.639
$ws.Range("D52").Value = -7.585000000000001
$ws.Range("C54").Value = -13.376
$ws.Range("A55").Value = -21.724
$ws.Range("B55").Value = 6.391999999999999
$ws.Range("D56").Value = -7.972
$ws.Range("A57").Value = -21.349
$ws.Range("C57").Value = -13.051
$ws.Range("C62").Value = -13.61
$ws.Range("C63").Value = -12.277
$ws.Range("C67").Value = -11.224
$ws.Range("A69").Value = -21.42899999999999
$ws.Range("B69").Value = 6.43
$ws.Range("B70").Value = 6.08
$ws.Range("C70").Value = -11.114
$ws.Range("A76").Value = -21.706
$ws.Range("B76").Value = 6.703
$ws.Range("A78").Value = -20.786
$ws.Range("B81").Value = 5.497999999999999
$ws.Range("C81").Value = -12.632
$ws.Range("A82").Value = -21.896
$ws.Range("A83").Value = -20.927
$ws.Range("B83").Value = 7.124
$ws.Range("D84").Value = -8.181000000000001
$ws.Range("C88").Value = -13.324
$ws.Range("A93").Value = -21.77
$ws.Range("B94").Value = 6.742
$ws.Range("C96").Value = -12.998
$ws.Range("A97").Value = -21.795
$ws.Range("C99").Value = -12.729
$ws.Range("B100").Value = 6.196
$ws.Range("D100").Value = -8.503
$ws.Range("D101").Value = -7.742
$ws.Range("B102").Value = 6.922
